$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price/volume data (and two row swaps: rows 32/33, rows 37/38)

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.497.68"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +1.93%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.865.38"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +2.57%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.44%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "316.30"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +2.34%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.003"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.47%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4680"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +2.10%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07382"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +2.21%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +3.48%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07934"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +1.38%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.899.55"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +7.99%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.426"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.599"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +2.05%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "92.66"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.99%  "
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.35%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008930"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.27%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.91"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +3.13%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "27.536.34"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +3.51%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.162"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.56%  "
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.46%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.119.59"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +2.55%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "153.49"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +1.09%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.876"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +1.43%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +2.46%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.092"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +1.38%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.173"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +1.31%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "117.00"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +1.51%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08909"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.37%  "
$ws.Range("B32").Value = "ImmutableX"
$ws.Range("C32").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7576"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +5.89%  "
$ws.Range("B33").Value = "HuobiToken"
$ws.Range("C33").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.027"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +2.66%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.167"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +3.40%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.488"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +1.70%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +10.82%  "
$ws.Range("B37").Value = "TrustWalletToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.083"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.57%  "
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01967"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +2.43%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05274"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +0.50%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +2.40%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.177"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.75%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5221"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +1.44%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1647"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +1.44%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.385"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +2.92%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4876"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +1.39%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.32"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +1.98%  "
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.49%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "103.93"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +1.22%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.657"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +2.50%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06267"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.22%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "65.99"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +2.44%  "
